$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6462
$ws.Range("C25").Value = 1007
$ws.Range("D25").Value = 6009416
$ws.Range("E25").Value = 929.9622407923243
$ws.Range("F25").Value = 9.692751655067045
$ws.Range("G25").Value = 7.356076759061825
$ws.Range("H25").Value = 25.85570299021163
